$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @'
questions = [
    {
        "title": "While parsing JSON from an API, you receive the names and ages of users. You need to send this data to another API as a single vector. However, you face a problem. The data comes in two types: u32 for the ages and String for the names.  What can you do to hold values of different types in the same vector?",
        "ques_type": 2,
        "options": [
            "Define a custom enum that can contain either u32 or String, and make that the type of the vector.",
            "Add the data using the .push() method.",
            "Use variable shadowing to create a vector that accepts both data types.",
            "Convert them both to u8 before adding them to the vector."
        ],
        "score": "Define a custom enum that can contain either u32 or String, and make that the type of the vector."
    },
    {
        "title": "You are writing a terminal application that only supports Bash (Bourne Shell), and you want to detect what shell the user is running. However, the response variable in the code below is not accessible when you try to print it. How can you fix the bug? use std::io::{stdin}\n \nfn main() {\n   let mut s = String::new()\n   println!(\"What shell do you use?\")\n   stdin().read_line(&ampmut s)\n \n   if s.trim().to_lowercase() != \"bash\" {\n       let response = String::from(\"Only bash is supported\")\n   } else {\n       let response = String::from(\"Great, your shell is supported!\")\n   }\n \n   println!(\"{}\", response)\n}",
        "ques_type": 2,
        "options": [
            "Make the response variable a constant so it can be used anywhere in the program.",
            "Switch to use match syntax instead of if for the conditional.",
            "Create the response variable as a mutable string before the if {} block, so it is in the scope of the main() function.",
            "Change the original string, s, to be immutable."
        ],
        "score": "Create the response variable as a mutable string before the if {} block, so it is in the scope of the main() function."
    },
    {
        "title": "You wrote some Rust code to scan a local network for suspicious activity. Other developers like your code and want to use it, so you decide to package it as a module. How can you use Cargo to create a new module for your code?",
        "ques_type": 2,
        "options": [
            "Add the code to crates.io by submitting a PR (pull request).",
            "Submit a PR with Git to add the code to crates.io.",
            "Run the cargo new command with the --bin option.",
            "Run the cargo new command with the --lib option."
        ],
        "score": "Run the cargo new command with the --lib option."
    },
    {
        "title": "You are developing an application for Windows that often requires unsafe blocks to access native C++ functions via the Foreign Function Interface (FFI). When you express concern, the product manager says not to worry, as the borrow checker will protect it. True or false: The borrow checker protects code wrapped in an unsafe block from memory bugs.",
        "ques_type": 11,
        "options": [
            "true",
            "false"
        ],
        "score": "False"
    }
]
'@

$ws.Range("A1").Value = $newText
$ws.Range("A1").Style = "Normal"
$ws.Rows(1).AutoFit()
$ws.Range("A2").ClearContents()
